$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 419894.5
$ws.Range("J3").Value = 419894.5
$ws.Range("L3").Value = 419894.5
$ws.Range("N3").Value = -420122.5
$ws.Range("H9").Value = 7549.875
$ws.Range("I9").Value = 8221.286
$ws.Range("K9").Value = 8221.286
$ws.Range("M9").Value = -8052.286
$ws.Range("H58").Value = 446.875
$ws.Range("I58").Value = 367.85715
$ws.Range("K58").Value = 1103.57145
$ws.Range("M58").Value = -953.5714499999999
$ws.Range("H64").Value = 10579.4
$ws.Range("I64").Value = 10579.4
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 10579.4
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -10331.4
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 10579.4
$ws.Range("I67").Value = 10579.4
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 10579.4
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -9721.4
$ws.Range("H98").Value = 1431.3529
$ws.Range("I98").Value = 1431.3529
$ws.Range("K98").Value = 1431.3529
$ws.Range("M98").Value = 66.64709999999991
$ws.Range("H101").Value = 1131.8
$ws.Range("I101").Value = 1189.875
$ws.Range("J101").Value = 899.5
$ws.Range("K101").Value = 3569.625
$ws.Range("L101").Value = 2698.5
$ws.Range("M101").Value = -1947.625
$ws.Range("N101").Value = -5942.5
$ws.Range("H102").Value = 419894.5
$ws.Range("J102").Value = 419894.5
$ws.Range("L102").Value = 419894.5
$ws.Range("N102").Value = -426384.5
$ws.Range("H103").Value = 805
$ws.Range("J103").Value = 805
$ws.Range("L103").Value = 2415
$ws.Range("N103").Value = -3587
$ws.Range("H113").Value = 1971.4286
$ws.Range("I113").Value = 1933.3334
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1933.3334
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 1320.6666
$ws.Range("N113").Value = -8508
$ws.Range("H122").Value = 1431.3529
$ws.Range("I122").Value = 1431.3529
$ws.Range("K122").Value = 4294.0587
$ws.Range("M122").Value = -1844.0587

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 63.333332
$ws.Range("I4").Value = 63.333332
$ws.Range("K4").Value = 63.333332
$ws.Range("M4").Value = 52.666668
$ws.Range("H32").Value = 4606.924
$ws.Range("I32").Value = 3704.1353
$ws.Range("K32").Value = 3704.1353
$ws.Range("M32").Value = -3417.1353
$ws.Range("H88").Value = 6772.8335
$ws.Range("I88").Value = 15330.714
$ws.Range("J88").Value = 1326.909
$ws.Range("K88").Value = 15330.714
$ws.Range("L88").Value = 1326.909
$ws.Range("M88").Value = -14924.714
$ws.Range("N88").Value = -2138.909
$ws.Range("H91").Value = 6772.8335
$ws.Range("I91").Value = 15330.714
$ws.Range("J91").Value = 1326.909
$ws.Range("K91").Value = 15330.714
$ws.Range("L91").Value = 1326.909
$ws.Range("M91").Value = -13926.714
$ws.Range("N91").Value = -4134.909
$ws.Range("H132").Value = 47697836
$ws.Range("I132").Value = 14043.944
$ws.Range("K132").Value = 42131.83199999999
$ws.Range("M132").Value = -39601.83199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2538.1428
$ws.Range("I20").Value = 2941.75
$ws.Range("J20").Value = 2000
$ws.Range("K20").Value = 2941.75
$ws.Range("L20").Value = 2000
$ws.Range("M20").Value = -2694.75
$ws.Range("N20").Value = -2494
$ws.Range("H22").Value = 545.125
$ws.Range("J22").Value = 838.5
$ws.Range("L22").Value = 838.5
$ws.Range("N22").Value = -1184.5
$ws.Range("H94").Value = 972.84
$ws.Range("I94").Value = 407.47058
$ws.Range("K94").Value = 407.47058
$ws.Range("M94").Value = 43.52942000000002
$ws.Range("H140").Value = 148233.4
$ws.Range("J140").Value = 148233.4
$ws.Range("L140").Value = 148233.4
$ws.Range("N140").Value = -158593.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 7069.8667
$ws.Range("I22").Value = 8042.154
$ws.Range("K22").Value = 8042.154
$ws.Range("M22").Value = -7692.154
$ws.Range("H31").Value = 4307.1934
$ws.Range("I31").Value = 3283.8823
$ws.Range("K31").Value = 3283.8823
$ws.Range("M31").Value = -2988.8823
$ws.Range("H34").Value = 4307.1934
$ws.Range("I34").Value = 3283.8823
$ws.Range("K34").Value = 3283.8823
$ws.Range("M34").Value = -3081.8823
$ws.Range("H62").Value = 3880
$ws.Range("I62").Value = 3166.6667
$ws.Range("K62").Value = 3166.6667
$ws.Range("M62").Value = -2542.6667
$ws.Range("H65").Value = 3880
$ws.Range("I65").Value = 3166.6667
$ws.Range("K65").Value = 15833.3335
$ws.Range("M65").Value = -12713.3335

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 292.4
$ws.Range("I8").Value = 292.4
$ws.Range("K8").Value = 877.1999999999999
$ws.Range("M8").Value = -738.1999999999999
$ws.Range("H18").Value = 1431.9166
$ws.Range("I18").Value = 806.625
$ws.Range("J18").Value = 2682.5
$ws.Range("K18").Value = 2419.875
$ws.Range("L18").Value = 8047.5
$ws.Range("M18").Value = -2250.875
$ws.Range("N18").Value = -8385.5
$ws.Range("H70").Value = 7499.25
$ws.Range("I70").Value = 99.5
$ws.Range("J70").Value = 14899
$ws.Range("K70").Value = 298.5
$ws.Range("L70").Value = 44697
$ws.Range("M70").Value = 16.5
$ws.Range("N70").Value = -45327
$ws.Range("H73").Value = 7499.25
$ws.Range("I73").Value = 99.5
$ws.Range("J73").Value = 14899
$ws.Range("K73").Value = 298.5
$ws.Range("L73").Value = 44697
$ws.Range("M73").Value = 793.5
$ws.Range("N73").Value = -46881
$ws.Range("H122").Value = 707.3333
$ws.Range("J122").Value = 788.75
$ws.Range("L122").Value = 7098.75
$ws.Range("N122").Value = -11998.75
$ws.Range("H131").Value = 31181.184
$ws.Range("I131").Value = 116385.445
$ws.Range("J131").Value = 4738.483
$ws.Range("K131").Value = 349156.335
$ws.Range("L131").Value = 14215.449
$ws.Range("M131").Value = -344116.335
$ws.Range("N131").Value = -24295.449

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 48333.668
$ws.Range("H50").Value = 48333.668
$ws.Range("H70").Value = 503642.5
$ws.Range("I70").Value = 669857
$ws.Range("J70").Value = 4999
$ws.Range("K70").Value = 669857
$ws.Range("L70").Value = 4999
$ws.Range("M70").Value = -669587
$ws.Range("N70").Value = -5539
$ws.Range("H73").Value = 503642.5
$ws.Range("I73").Value = 669857
$ws.Range("J73").Value = 4999
$ws.Range("K73").Value = 669857
$ws.Range("L73").Value = 4999
$ws.Range("M73").Value = -668921
$ws.Range("N73").Value = -6871
$ws.Range("H97").Value = 1789.1892
$ws.Range("I97").Value = 1307.28
$ws.Range("K97").Value = 1307.28
$ws.Range("M97").Value = -811.28
$ws.Range("H122").Value = 2839.6924
$ws.Range("I122").Value = 2515.75
$ws.Range("K122").Value = 7547.25
$ws.Range("M122").Value = -5097.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4600
$ws.Range("I40").Value = 4600
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 4600
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -4464
$ws.Range("N40").ClearContents()
$ws.Range("H46").Value = 1369.5555
$ws.Range("J46").Value = 3719
$ws.Range("L46").Value = 3719
$ws.Range("N46").Value = -4095
$ws.Range("H68").Value = 2200
$ws.Range("J68").Value = 2200
$ws.Range("L68").Value = 2200
$ws.Range("N68").Value = -3698
$ws.Range("H71").Value = 2200
$ws.Range("J71").Value = 2200
$ws.Range("L71").Value = 11000
$ws.Range("N71").Value = -18488
$ws.Range("H93").Value = 1870.25
$ws.Range("I93").Value = 1098
$ws.Range("K93").Value = 1098
$ws.Range("M93").Value = 150
$ws.Range("H122").Value = 4863.647
$ws.Range("I122").Value = 4081
$ws.Range("K122").Value = 12243
$ws.Range("M122").Value = -9793
$ws.Range("H136").Value = 11332.333
$ws.Range("I136").Value = 8998
$ws.Range("K136").Value = 26994
$ws.Range("M136").Value = -24444

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6699.9
$ws.Range("I62").Value = 4500
$ws.Range("J62").Value = 7249.875
$ws.Range("K62").Value = 4500
$ws.Range("L62").Value = 7249.875
$ws.Range("M62").Value = -3876
$ws.Range("N62").Value = -8497.875
$ws.Range("H65").Value = 6699.9
$ws.Range("I65").Value = 4500
$ws.Range("J65").Value = 7249.875
$ws.Range("K65").Value = 22500
$ws.Range("L65").Value = 36249.375
$ws.Range("M65").Value = -19380
$ws.Range("N65").Value = -42489.375
$ws.Range("H81").Value = 1441.4667
$ws.Range("J81").Value = 2034.6
$ws.Range("L81").Value = 4069.2
$ws.Range("N81").Value = -6191.2
$ws.Range("H84").Value = 1441.4667
$ws.Range("J84").Value = 2034.6
$ws.Range("L84").Value = 20346
$ws.Range("N84").Value = -30954
$ws.Range("H132").Value = 4074.8262
$ws.Range("I132").Value = 6058.5
$ws.Range("K132").Value = 18175.5
$ws.Range("M132").Value = -15645.5
$ws.Range("H136").Value = 3454.8333
$ws.Range("I136").Value = 987.5333000000001
$ws.Range("K136").Value = 2962.5999
$ws.Range("M136").Value = -412.5999000000002

Write-Host "Applied all Famfrit_Profits updates"
